$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C26 changes style (number format) from "YYYY-MM-DD" (s=3) to
# "YYYY-MM-DD HH:MM:SS" (s=2), value unchanged.
$ws.Range("C26").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 27
$ws.Range("A27").Value = 770.419
$ws.Range("B27").Value = 692.068
$ws.Range("C27").Value = 45758
$ws.Range("C27").NumberFormat = "YYYY-MM-DD"
